$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Referensi")

# --- Row 24: derHugo / Unity 2019 blender->FBX import bug workaround ---
$ws.Range("A24").Value = "derHugo"
$ws.Range("B24").Value = "https://stackoverflow.com/questions/55752495/unity-2019-1-0f2-blender-could-not-convert-the-blend-file-to-fbx-file"
$ws.Range("C24").Value = "Unity 2019 import blender into FBX bug workaround"
$ws.Range("D24").Value = "Unity 2019 terdapat masalah dengan impor Blender sebelum versi 2.8"
$ws.Hyperlinks.Add($ws.Range("B24"), "https://stackoverflow.com/questions/55752495/unity-2019-1-0f2-blender-could-not-convert-the-blend-file-to-fbx-file")
$ws.Range("B24").Style = $ws.Range("B23").Style

# --- Row 25: Manash Kumar Mandal / Arduino serial communication C++ ---
$ws.Range("A25").Value = "Manash Kumar Mandal"
$ws.Range("B25").Value = "https://blog.manash.me/serial-communication-with-an-arduino-using-c-on-windows-d08710186498"
$ws.Range("C25").Value = "Arduino Serial communication C++"
$ws.Range("D25").Value = "Komunikasi arduino menggunakan C++ melalui serial, buat Godot C++"
$ws.Hyperlinks.Add($ws.Range("B25"), "https://blog.manash.me/serial-communication-with-an-arduino-using-c-on-windows-d08710186498")
$ws.Range("B25").Style = $ws.Range("B23").Style

# --- Column width tweaks ---
$ws.Columns.Item(1).ColumnWidth = 30.85
$ws.Columns.Item(4).ColumnWidth = 58.95

# --- View / selection state ---
$ws.Range("D21").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
